$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- Cells that change FROM a text placeholder ("0" / "***.*") TO a real number: set value + number format ---
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 3).NumberFormat = '#,##0'
$ws.Cells.Item(15, 9).Value = 1
$ws.Cells.Item(15, 9).NumberFormat = '#,##0'
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 12).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 9).NumberFormat = '#,##0'
$ws.Cells.Item(16, 12).Value = -50
$ws.Cells.Item(16, 12).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 4).NumberFormat = '#,##0'
$ws.Cells.Item(18, 5).Value = 300
$ws.Cells.Item(18, 5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(18, 10).Value = 1
$ws.Cells.Item(18, 10).NumberFormat = '#,##0'
$ws.Cells.Item(18, 11).Value = 400
$ws.Cells.Item(18, 11).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(18, 12).Value = 400
$ws.Cells.Item(18, 12).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(20, 9).Value = 2
$ws.Cells.Item(20, 9).NumberFormat = '#,##0'
$ws.Cells.Item(20, 10).Value = 10
$ws.Cells.Item(20, 10).NumberFormat = '#,##0'
$ws.Cells.Item(20, 11).Value = -80
$ws.Cells.Item(20, 11).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(23, 12).Value = -100
$ws.Cells.Item(23, 12).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 3).NumberFormat = '#,##0'
$ws.Cells.Item(27, 9).Value = 1
$ws.Cells.Item(27, 9).NumberFormat = '#,##0'
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 12).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 3).NumberFormat = '#,##0'
$ws.Cells.Item(28, 9).Value = 1
$ws.Cells.Item(28, 9).NumberFormat = '#,##0'
$ws.Cells.Item(28, 10).Value = 1
$ws.Cells.Item(28, 10).NumberFormat = '#,##0'
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 11).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(28, 12).Value = -66.666666666666
$ws.Cells.Item(28, 12).NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Plain value updates (style/number format unchanged) ---
$ws.Cells.Item(15, 6).Value = 4
$ws.Cells.Item(15, 8).Value = 300
$ws.Cells.Item(15, 14).Value = -85.714285714285
$ws.Cells.Item(16, 4).Value = 4
$ws.Cells.Item(16, 5).Value = -75
$ws.Cells.Item(16, 6).Value = 3
$ws.Cells.Item(16, 7).Value = 11
$ws.Cells.Item(16, 8).Value = -72.727272727272
$ws.Cells.Item(16, 10).Value = 5
$ws.Cells.Item(16, 11).Value = -80
$ws.Cells.Item(16, 13).Value = -94.117647058823
$ws.Cells.Item(16, 14).Value = -96.875
$ws.Cells.Item(17, 3).Value = 4
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(17, 5).Value = 100
$ws.Cells.Item(17, 6).Value = 23
$ws.Cells.Item(17, 7).Value = 21
$ws.Cells.Item(17, 8).Value = 9.523809523809
$ws.Cells.Item(17, 9).Value = 8
$ws.Cells.Item(17, 10).Value = 8
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = -38.461538461538
$ws.Cells.Item(17, 13).Value = 14.285714285714
$ws.Cells.Item(17, 14).Value = -61.904761904761
$ws.Cells.Item(18, 3).Value = 4
$ws.Cells.Item(18, 6).Value = 9
$ws.Cells.Item(18, 7).Value = 5
$ws.Cells.Item(18, 8).Value = 80
$ws.Cells.Item(18, 9).Value = 5
$ws.Cells.Item(18, 13).Value = -58.333333333333
$ws.Cells.Item(18, 14).Value = -87.179487179487
$ws.Cells.Item(19, 3).Value = 10
$ws.Cells.Item(19, 4).Value = 9
$ws.Cells.Item(19, 5).Value = 11.111111111111
$ws.Cells.Item(19, 6).Value = 40
$ws.Cells.Item(19, 7).Value = 25
$ws.Cells.Item(19, 8).Value = 60
$ws.Cells.Item(19, 9).Value = 17
$ws.Cells.Item(19, 10).Value = 11
$ws.Cells.Item(19, 11).Value = 54.545454545454
$ws.Cells.Item(19, 12).Value = 183.333333333333
$ws.Cells.Item(19, 13).Value = 41.666666666666
$ws.Cells.Item(19, 14).Value = -72.131147540983
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 10
$ws.Cells.Item(20, 5).Value = -80
$ws.Cells.Item(20, 6).Value = 12
$ws.Cells.Item(20, 7).Value = 20
$ws.Cells.Item(20, 8).Value = -40
$ws.Cells.Item(20, 12).Value = -80
$ws.Cells.Item(20, 13).Value = -85.714285714285
$ws.Cells.Item(20, 14).Value = -96.296296296296
$ws.Cells.Item(21, 3).Value = 22
$ws.Cells.Item(21, 4).Value = 26
$ws.Cells.Item(21, 5).Value = -15.384615384615
$ws.Cells.Item(21, 6).Value = 91
$ws.Cells.Item(21, 8).Value = 9.638554216867
$ws.Cells.Item(21, 9).Value = 34
$ws.Cells.Item(21, 10).Value = 35
$ws.Cells.Item(21, 11).Value = -2.857142857142
$ws.Cells.Item(21, 12).Value = 3.030303030303
$ws.Cells.Item(21, 13).Value = -45.161290322580
$ws.Cells.Item(21, 14).Value = -84.259259259259
$ws.Cells.Item(24, 4).Value = 13
$ws.Cells.Item(24, 5).Value = 92.307692307692
$ws.Cells.Item(24, 6).Value = 72
$ws.Cells.Item(24, 7).Value = 67
$ws.Cells.Item(24, 8).Value = 7.462686567164
$ws.Cells.Item(24, 9).Value = 30
$ws.Cells.Item(24, 10).Value = 17
$ws.Cells.Item(24, 11).Value = 76.470588235294
$ws.Cells.Item(24, 12).Value = -28.571428571428
$ws.Cells.Item(24, 13).Value = 7.142857142857
$ws.Cells.Item(25, 3).Value = 5
$ws.Cells.Item(25, 4).Value = 4
$ws.Cells.Item(25, 5).Value = 25
$ws.Cells.Item(25, 6).Value = 15
$ws.Cells.Item(25, 7).Value = 14
$ws.Cells.Item(25, 8).Value = 7.142857142857
$ws.Cells.Item(25, 9).Value = 6
$ws.Cells.Item(25, 10).Value = 5
$ws.Cells.Item(25, 11).Value = 20
$ws.Cells.Item(25, 12).Value = -72.727272727272
$ws.Cells.Item(26, 3).Value = 12
$ws.Cells.Item(26, 4).Value = 13
$ws.Cells.Item(26, 5).Value = -7.692307692307
$ws.Cells.Item(26, 6).Value = 47
$ws.Cells.Item(26, 8).Value = -11.320754716981
$ws.Cells.Item(26, 9).Value = 20
$ws.Cells.Item(26, 10).Value = 24
$ws.Cells.Item(26, 11).Value = -16.666666666666
$ws.Cells.Item(26, 12).Value = 53.846153846153
$ws.Cells.Item(26, 13).Value = 17.647058823529
$ws.Cells.Item(27, 6).Value = 4
$ws.Cells.Item(27, 8).Value = 300
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 7).Value = 7
$ws.Cells.Item(28, 8).Value = -71.428571428571
$ws.Cells.Item(44, 10).Value = 459
$ws.Cells.Item(44, 11).Value = -46.315789473684
$ws.Cells.Item(44, 12).Value = -56.244041944709
$ws.Cells.Item(44, 13).Value = -88.473129080863
$ws.Cells.Item(44, 14).Value = -85.073170731707
$ws.Cells.Item(46, 10).Value = 1300
$ws.Cells.Item(46, 11).Value = -50.943396226415
$ws.Cells.Item(46, 12).Value = -64.078474716772
$ws.Cells.Item(46, 13).Value = -85.582788066984
$ws.Cells.Item(46, 14).Value = -85.900216919739
